$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value2 = 130654938
$ws.Range("B9").Value2 = 79243
$ws.Range("E9").Value2 = 6425
$ws.Range("F9").Value2 = 'Garnlav'
$ws.Range("G9").Value2 = 'Alectoria sarmentosa'
$ws.Range("H9").Value2 = '(Ach.) Ach.'
$ws.Range("J9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("Q9").Value2 = 440117
$ws.Range("R9").Value2 = 7053967
$ws.Range("AC9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AM9").Value2 = 'Gren på levande träd'
$ws.Range("AO9").Value2 = 'Branch on living tree # Picea abies'

# Row 10
$ws.Range("A10").Value2 = 130654930
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("M10").Value2 = 'färska spår'
$ws.Range("N10").ClearContents()
$ws.Range("Q10").Value2 = 440168
$ws.Range("R10").Value2 = 7053746
$ws.Range("AC10").Value2 = 'Ringhack (savhack), enstaka färska, några meter upp på en gran vid en hyggeskant.'
$ws.Range("AH10").Value2 = 'Granskog'
$ws.Range("AJ10").Value2 = 'gran'
$ws.Range("AK10").Value2 = 'Picea abies'
$ws.Range("AM10").Value2 = 'Trädstam på levande träd'
$ws.Range("AO10").Value2 = 'Stem on living tree # Picea abies'
$ws.Range("AW10").Value2 = 'Kristian Zackrisson'
$ws.Range("AX10").Value2 = 'Kristian Zackrisson'

# Row 11
$ws.Range("A11").Value2 = 130654517
$ws.Range("B11").Value2 = 57884
$ws.Range("E11").Value2 = 100109
$ws.Range("F11").Value2 = 'Tretåig hackspett'
$ws.Range("G11").Value2 = 'Picoides tridactylus'
$ws.Range("H11").Value2 = '(Linnaeus, 1758)'
$ws.Range("J11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("Q11").Value2 = 440178
$ws.Range("R11").Value2 = 7053979
$ws.Range("AC11").Value2 = 'Ringhack äldre'
$ws.Range("AF11").ClearContents()
$ws.Range("AH11").ClearContents()
$ws.Range("AJ11").ClearContents()
$ws.Range("AK11").ClearContents()
$ws.Range("AM11").ClearContents()
$ws.Range("AO11").ClearContents()
$ws.Range("AW11").Value2 = 'Benny Öwre'
$ws.Range("AX11").Value2 = 'Benny Öwre'

# Row 17
$ws.Range("A17").Value2 = 130654518
$ws.Range("B17").Value2 = 57884
$ws.Range("E17").Value2 = 100109
$ws.Range("F17").Value2 = 'Tretåig hackspett'
$ws.Range("G17").Value2 = 'Picoides tridactylus'
$ws.Range("H17").Value2 = '(Linnaeus, 1758)'
$ws.Range("J17").ClearContents()
$ws.Range("K17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("Q17").Value2 = 440177
$ws.Range("R17").Value2 = 7054022
$ws.Range("AC17").Value2 = 'Ringhack äldre'
$ws.Range("AF17").ClearContents()
$ws.Range("AH17").ClearContents()
$ws.Range("AJ17").ClearContents()
$ws.Range("AK17").ClearContents()
$ws.Range("AM17").ClearContents()
$ws.Range("AO17").ClearContents()
$ws.Range("AW17").Value2 = 'Benny Öwre'
$ws.Range("AX17").Value2 = 'Benny Öwre'

# Row 18
$ws.Range("A18").Value2 = 130654935
$ws.Range("B18").Value2 = 79243
$ws.Range("E18").Value2 = 6425
$ws.Range("F18").Value2 = 'Garnlav'
$ws.Range("G18").Value2 = 'Alectoria sarmentosa'
$ws.Range("H18").Value2 = '(Ach.) Ach.'
$ws.Range("J18").ClearContents()
$ws.Range("K18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("Q18").Value2 = 439862
$ws.Range("R18").Value2 = 7054226
$ws.Range("AC18").Value2 = 'Enstaka bålar på gran.'
$ws.Range("AF18").ClearContents()
$ws.Range("AH18").Value2 = 'Granskog'
$ws.Range("AJ18").Value2 = 'gran'
$ws.Range("AK18").Value2 = 'Picea abies'
$ws.Range("AM18").Value2 = 'Gren på levande träd'
$ws.Range("AO18").Value2 = 'Branch on living tree # Picea abies'
$ws.Range("AW18").Value2 = 'Kristian Zackrisson'
$ws.Range("AX18").Value2 = 'Kristian Zackrisson'

# Row 22
$ws.Range("A22").Value2 = 130654933
$ws.Range("B22").Value2 = 57881
$ws.Range("E22").Value2 = 100049
$ws.Range("F22").Value2 = 'Spillkråka'
$ws.Range("G22").Value2 = 'Dryocopus martius'
$ws.Range("M22").Value2 = 'äldre spår'
$ws.Range("Q22").Value2 = 439962
$ws.Range("R22").Value2 = 7053998
$ws.Range("AC22").Value2 = 'Äldre hackspår i stambasen av en stående död gran med full längd.'
$ws.Range("AM22").Value2 = 'Stående död trädstam/högstubbe'
$ws.Range("AO22").Value2 = 'Standing dead tree/snags # Picea abies'

# Row 23
$ws.Range("A23").Value2 = 130654931
$ws.Range("B23").Value2 = 57884
$ws.Range("E23").Value2 = 100109
$ws.Range("F23").Value2 = 'Tretåig hackspett'
$ws.Range("G23").Value2 = 'Picoides tridactylus'
$ws.Range("M23").Value2 = 'färska spår'
$ws.Range("Q23").Value2 = 440195
$ws.Range("R23").Value2 = 7053721
$ws.Range("AC23").Value2 = 'Ringhack (savhack), färska, enstaka på en gran vid en hyggeskant.'
$ws.Range("AM23").Value2 = 'Trädstam på levande träd'
$ws.Range("AO23").Value2 = 'Stem on living tree # Picea abies'
